# Repull data, push all data, mean calculation
# Updates the "dSF" column (column F) values for a set of rows to reflect
# the repulled/updated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$updates = @{
    18 = -3
    20 = 1
    33 = -5
    36 = 10
    40 = -5
    42 = -1
    49 = -3
    51 = -6
    54 = -6
    56 = 8
    57 = 3
    60 = 5
    62 = -1
    63 = -6
    64 = -6
    66 = 10
    69 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
